$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, shifting the existing rows 145-150 down to 146-151
$ws.Rows("145:145").Insert()

# Copy the style (number format) of the date cell from the row below (now D146,
# which held the original D145 style) onto the new D145 cell.
$ws.Range("D146").Copy()
$ws.Range("D145").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the constant columns shared by every row of this dataset
$ws.Range("A145").Value = 7
$ws.Range("B145").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C145").Value = "Ñuble"
$ws.Range("D145").Value = 44931
$ws.Range("E145").Value = 16
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100109
$ws.Range("H145").Value = "Uva"
$ws.Range("I145").Value = 100109001
$ws.Range("J145").Value = "Uva"
$ws.Range("K145").Value = "Flame Seedless"
$ws.Range("L145").Value = "Primera"
$ws.Range("M145").Value = 100
$ws.Range("N145").Value = 15000
$ws.Range("O145").Value = 16000
$ws.Range("P145").Value = 15500
$ws.Range("Q145").Value = "`$/bandeja 10 kilos"
$ws.Range("R145").Value = "Provincia de Limarí"
$ws.Range("S145").Value = 1550
$ws.Range("T145").Value = 10
